$d = $word.ActiveDocument

# The three inline pictures in the headers/footers (Pearson logo x2, BTEC
# logo x1) need their shape "name" (docPr/cNvPr @name) swapped:
#   image2.png -> image1.png   (both Pearson logo instances)
#   image1.jpg -> image2.jpg   (BTEC logo)
# Word's object model does not expose a settable Name on InlineShape
# (that only exists on the floating Shape type), so the rename is done by
# editing the underlying part XML via Document.WordOpenXML, which is a
# faithful flat-OPC serialization of the whole package.

$xml = $d.WordOpenXML

$xml = $xml -replace 'name="image2\.png"', 'name="image1.png"'
$xml = $xml -replace 'name="image1\.jpg"', 'name="image2.jpg"'

$d.WordOpenXML = $xml
